# "Small tweeks to the formating" -- rename sheets, tidy up the
# Competitors summary sheet (fonts/row heights/alignment), set a
# print area, tighten margins, and lock the summary sheet down.

$wb = $excel.ActiveWorkbook

$wsBracket = $wb.Worksheets.Item(1)
$wsCompetitors = $wb.Worksheets.Item(2)

# Rename the sheets (formulas referencing the old names update automatically)
$wsBracket.Name = "Bracket"
$wsCompetitors.Name = "Competitors"

# Print area for the Competitors sheet
$wsCompetitors.PageSetup.PrintArea = '$A:$D'

# Competitors sheet formatting tweaks
# Row 2 and 3 shrink from 18pt "Arial Unicode MS" down to a tidier 14pt Calibri (Body)
$wsCompetitors.Range("A2:D3").Font.Size = 14
$wsCompetitors.Range("A2:D3").Font.Name = "Calibri (Body)"

# Rows get shorter now that the font is smaller
$wsCompetitors.Range("A2").EntireRow.RowHeight = 19
$wsCompetitors.Range("A3").EntireRow.RowHeight = 19

# B1 (the mat/bracket number) becomes left-aligned
$wsCompetitors.Range("B1").HorizontalAlignment = -4131

# Tighten up the page margins on the Competitors sheet
$wsCompetitors.PageSetup.LeftMargin = 14.4
$wsCompetitors.PageSetup.RightMargin = 14.4
$wsCompetitors.PageSetup.TopMargin = 18

# Lock the Competitors sheet down like the Bracket sheet already is
$wsCompetitors.Protect("CFAF")
